$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = -0.2920432377888471
$ws.Range("J2").Value = 0.1126628593094076
$ws.Range("K2").Value = -0.6791700562381042
$ws.Range("L2").Value = 2.591575756073917

# Row 12
$ws.Range("I12").Value = -0.09340982480822607
$ws.Range("J12").Value = 0.03509831356904334
$ws.Range("K12").Value = -0.3797560568647471
$ws.Range("L12").Value = 1.922929931946162

# Row 13
$ws.Range("I13").Value = -0.1816609049986487
$ws.Range("J13").Value = 0.07681378770212414
$ws.Range("K13").Value = -0.6949592617086884
$ws.Range("L13").Value = 2.366446696117301

# Row 15
$ws.Range("I15").Value = -0.2108672473204946
$ws.Range("J15").Value = 0.07318313726210456
$ws.Range("K15").Value = -0.4092894075278155
$ws.Range("L15").Value = 2.10392564033197

# Row 18
$ws.Range("I18").Value = -0.4762472372781644
$ws.Range("J18").Value = 0.1369078737662411
$ws.Range("K18").Value = 0.00321483838635471
$ws.Range("L18").Value = 1.875983085926304
